# Move the test checklist from categorytest.txt into the Shapes Lab
# deck's speaker notes, attached to slide 2 ("Shapes Lab" instructions
# slide), per the commit "Move test in categorytest.txt into
# ShapesLab.pptx".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$notes = $s.NotesPage

$notesBody = $notes.Shapes.Placeholders.Item(1)

$lines = @(
    "Do the following from 3 environments:",
    "1. Entirely new environment",
    "2. Contains old shape info",
    "",
    "1. Add Category",
    "`t1.1 Add category then add a new shape inside, check the default name, close and open to check",
    "`t1.2 Add category then add a new shape inside, check the default name, set the category as default, close and open to check",
    "2. Rename Category",
    "`t2.1 Rename to a valid name",
    "`t2.2 Rename to an invaid name",
    "3. Import Category",
    "`t3.1 Import single category with namebox, no conflict name",
    "`t3.2 Import multiple categories with nameboxes, no confilct name",
    "`t3.3 Import single category with namebox, conflict with one existed name",
    "`t3.4 Import multiple categories with nameboxes, some of the boxes conflict with existed name",
    "`t3.5 Import single category without namebox",
    "`t3.6 Import single category without namebox, conflict with exist name",
    "`t3.7 Import multiple categories without namebox",
    "4. Migrate Category"
)

$notesBody.TextFrame.TextRange.Text = [string]::Join("`n", $lines)
